# Fix formatting issues introduced when the "Importe" amounts were scraped:
# numbers were captured using an Argentina/Spain-style thousands/decimal
# separator ("1.234,56") but must be normalized to plain "1234.56" text.
# Also fix a few Razon social entries whose punctuation was scraped wrong.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H ("Importe"): normalize "1.234,56" -> "1234.56" -------------
# Drop the "." thousands separators, then turn the "," decimal separator
# into a ".". Re-apply the result as TEXT (not a number) so the cell keeps
# its original General/no-style formatting and the literal string content
# (e.g. trailing zeros) is preserved exactly.
for ($row = 2; $row -le 226; $row++) {
  $cell = $ws.Cells.Item($row, 8)
  $old = $cell.Text
  $new = $old.Replace(".", "").Replace(",", ".")
  if ($new -ne $old) {
    $cell.NumberFormat = "@"
    $cell.Value = $new
    $cell.Style = "Normal"
  }
}

# --- Column E ("Razon social"): punctuation fixes -------------------------
$ws.Cells.Item(64, 5).Value = "URUMAT SOCIEDAD SIMPLE DE BONASEGLA CATALINA. BONASEGLA LUCIANA Y BONASEGLA SILVIO"
$ws.Cells.Item(150, 5).Value = "RICCOTTI. MARIANA EDITH"
$ws.Cells.Item(168, 5).Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
